# Trade #55 closed at 2026-02-17 08:42:37 - unknown UNKNOWN +0.000%
#
# Appends the newly-closed trade (#55) to the "All Trades" and
# "MarketMaking" logs, and rolls the updated aggregate stats into the
# "Summary" and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet - roll the new trade into the top-level metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.56   # Current Capital
$summary.Range("B4").Value = -0.43     # Total P&L $
$summary.Range("B5").Value = -0.16     # Total P&L %
$summary.Range("B6").Value = 55        # Total Trades
$summary.Range("B7").Value = 21        # Winning Trades
$summary.Range("B9").Value = 38.18     # Win Rate %

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.56      # Capital
$status.Range("D4").Value = 55         # Trades
$status.Range("E4").Value = -0.43      # P&L $
$status.Range("F4").Value = -0.44      # P&L %
$status.Range("G4").Value = 38.18      # Win Rate %

# ---------------------------------------------------------------
# 3) Append the new trade row (#55 -> sheet row 56) to both the
#    "All Trades" log and the strategy-specific "MarketMaking" log.
# ---------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 55
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "'08:42:30"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.42
    $ws.Cells.Item($row, 7).Value = 0.47
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 11.9048
    $ws.Cells.Item($row, 10).Value = 0.05
    $ws.Cells.Item($row, 11).Value = 99.56
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 56

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 56
